# Audio volume reduced and phrases break fixed.
#
# The gaze/animate "prosody" template strings had a <break time="1s"/>
# inserted after every single <Gaze(...)>/<ANIMATE(...)> tag. The fix
# collapses those into a single <break time='1s'/> placed after the last
# tag (also switching the break tag to single-quoted attribute syntax),
# and re-aligns the PT/EN GAZE_* rows so each label's phrase actually
# matches its name (GAZE_PB = 2 gazes, GAZE_PBP = 3 gazes, GAZE_PBPB = 4
# gazes).

$wb = $excel.ActiveWorkbook

$utter = $wb.Worksheets.Item("Utterances")
$dropdowns = $wb.Worksheets.Item("Dropdowns Content")

# ---- New canonical phrase templates ------------------------------------
$gazeAnimate = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'> <Gaze(person3)><ANIMATE(|animation|)><break time='1s'/>  </prosody></prosody></prosody>"
$gazePB      = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'> <Gaze(|animation_person|)><Gaze(|animation_side|)><break time='1s'/>  </prosody></prosody></prosody>"
$gazePBP     = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'> <Gaze(|animation_person|)><Gaze(|animation_side|)><Gaze(|animation_person|)><break time='1s'/>  </prosody></prosody></prosody>"
$gazePBPalt  = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'> <Gaze(|animation_person|)><Gaze(|animation_side|)><Gaze(|animation_person|)><break time='1s'/> </prosody></prosody></prosody>"
$gazePBPB    = "<prosody pitch='|pitch|'><prosody rate='|rate|'><prosody volume='|volume|'> <Gaze(|animation_person|)><Gaze(|animation_side|)><Gaze(|animation_person|)><Gaze(|animation_side|)><break time='1s'/>  </prosody></prosody></prosody>"

# ---- Dropdown labels (B8:B13 on "Dropdowns Content") -------------------
$dropdowns.Range("B8").Value  = "PT:GAZE_PB"
$dropdowns.Range("B9").Value  = "EN:GAZE_PB"
$dropdowns.Range("B10").Value = "PT:GAZE_PBP"
$dropdowns.Range("B11").Value = "EN:GAZE_PBP"
$dropdowns.Range("B12").Value = "PT:GAZE_PBPB"
$dropdowns.Range("B13").Value = "EN:GAZE_PBPB"

# ---- Utterances sheet ---------------------------------------------------
# Rows 58-59: PT/EN:ANIMATION phrases (break moved to the end)
$utter.Range("D58").Value = $gazeAnimate
$utter.Range("D59").Value = $gazeAnimate

# Rows 60-61: GAZE_PB (2 gazes: person + side)
$utter.Range("B60").Value = "PT:GAZE_PB"
$utter.Range("D60").Value = $gazePB
$utter.Range("B61").Value = "EN:GAZE_PB"
$utter.Range("D61").Value = $gazePB

# Rows 62-63: GAZE_PBP (3 gazes: person + side + person)
$utter.Range("B62").Value = "PT:GAZE_PBP"
$utter.Range("D62").Value = $gazePBP
$utter.Range("B63").Value = "EN:GAZE_PBP"
$utter.Range("D63").Value = $gazePBPalt

# Rows 64-65: GAZE_PBPB (4 gazes: person + side + person + side)
$utter.Range("B64").Value = "PT:GAZE_PBPB"
$utter.Range("D64").Value = $gazePBPB
$utter.Range("B65").Value = "EN:GAZE_PBPB"
$utter.Range("D65").Value = $gazePBPB

$utter.Range("D59").Select()
